$d = $word.ActiveDocument

# 1. Date of Meeting: "1" + "3" + "/05/2024" -> "13/05/2024"
$d.Content.Find.Execute("Date of Meeting:1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Date of Meeting:13", 2)

# 2. " -> Dajana Kokomani, Xhulia Myftaraj" keep text the same (only proofErr split added around Kokomani,
#    no visible text change) - skip text replace.

# 3. "Hasanaj, Brend" text unchanged - only run split removed; skip text replace.

# 4. "Danja Korreshi, Aurel Kulemani" -> "Danja Korreshi, Aurel Kulemani, Olga Kolaj"
$d.Content.Find.Execute("Danja Korreshi, Aurel Kulemani", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Danja Korreshi, Aurel Kulemani, Olga Kolaj", 2)

# 5. "17/05/2024" stays same text, but merge runs "1"+"7"+"/05/2024..." -> single run.
$d.Content.Find.Execute("10:30 in the cafeteria.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "10:30 in the cafeteria.", 2)
